$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 6 into a brand new row 7 (same styles: A7 picks up the
# style index that A6 currently has).
$ws.Rows.Item(6).Copy()
$ws.Rows.Item(7).Insert(-4121)  # xlShiftDown

# Now re-point A6 at the same style used by the earlier date cells
# (A2:A5), which is style index 1 in the saved file.
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats

# Fill in row 7's actual data.
$ws.Range("A7").Value = 46068
$ws.Range("B7").Value = "Counting bits"
$ws.Range("C7").Value = "https://leetcode.com/problems/counting-bits/"

$ws.Hyperlinks.Add($ws.Range("C7"), "https://leetcode.com/problems/counting-bits/", "", "", "https://leetcode.com/problems/counting-bits/") | Out-Null
